$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1 ("Login Page Content"): move the selection, it is no longer
# the active/selected tab after this edit ---
$ws1.Range("A10").Select()

# --- Sheet 2 ("Form Input" -> "Login Functionality") ---
$ws2.Name = "Login Functionality"

# Remove the old rows 4-10 (mobile number / address rows no longer needed)
$ws2.Range("A4:B10").EntireRow.Delete()

# Replace the remaining data rows with the new login credential rows
$ws2.Range("A2").Value = "User Name Mail"
$ws2.Range("B2").Value = "test12312122@gmail.com"
$ws2.Range("A3").Value = "User password"
$ws2.Range("B3").Value = "Test@12345"

# Restore a basic page setup on the sheet
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Make "Login Functionality" the active/selected sheet with its new selection
$ws2.Activate()
$ws2.Range("A3").Select()
